$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 29 (shifts existing rows 29:33 down to 30:34,
# carrying their values/styles with them, matching the target diff where
# the former row 29 re-appears at row 30, etc.)
$ws.Rows(29).Insert()

# Populate the newly inserted row with the new entry (Fellype). Written in
# E, B, A, C, D order so any newly-registered shared strings land at the
# same indices the original authoring session produced.
$ws.Range("E29").Value = "kvlslck9"
$ws.Range("B29").Value = "Fellype"
$ws.Range("A29").Value = "61849ec5bee39a0026d534e7"
$ws.Range("C29").Value = 34957735
$ws.Range("D29").Value = 28

# The "posicao" (rank) column is a simple sequential count - renumber the
# rows that shifted down so the ranks stay contiguous (29..33 -> 30..34).
for ($r = 30; $r -le 34; $r++) {
    $ws.Cells.Item($r, 4).Value = $r - 1
}

# Match the saved selection state from the diff.
$ws.Range("D29").Select()
